# Update scripts with new TPM values (re-run of the NATMI LR-pairs export).
# The recomputed Tnfsf8-Tnfrsf8 table now includes an "ECs" sending cluster
# (rows 2-3) in addition to the original FAPs/MuSCs/Resolving-Mac senders
# (now rows 4-7, with refreshed expression numbers), plus a new
# "Resolving-Mac" sending cluster block (rows 8-9). Every numeric column is
# rewritten with the freshly computed values; the sending/target cluster
# labels shift down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnfsf8"
$ws.Range("C2").Value = "Tnfrsf8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.075119
$ws.Range("H2").Value = 0.225357
$ws.Range("I2").Value = 0.1238117789774295
$ws.Range("J2").Value = 0.1238117789774295
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.288105
$ws.Range("N2").Value = 0.8643149999999999
$ws.Range("O2").Value = 0.7917647996438357
$ws.Range("P2").Value = 0.7917647996438356
$ws.Range("Q2").Value = 0.021642159495
$ws.Range("R2").Value = 0.194779435455
$ws.Range("S2").Value = 0.0980298083756113
$ws.Range("T2").Value = 0.09802980837561129

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tnfsf8"
$ws.Range("C3").Value = "Tnfrsf8"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.075119
$ws.Range("H3").Value = 0.225357
$ws.Range("I3").Value = 0.1238117789774295
$ws.Range("J3").Value = 0.1238117789774295
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.07577199999999999
$ws.Range("N3").Value = 0.227316
$ws.Range("O3").Value = 0.2082352003561643
$ws.Range("P3").Value = 0.2082352003561643
$ws.Range("Q3").Value = 0.005691916868
$ws.Range("R3").Value = 0.051227251812
$ws.Range("S3").Value = 0.02578197060181815
$ws.Range("T3").Value = 0.02578197060181815

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Tnfsf8"
$ws.Range("C4").Value = "Tnfrsf8"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1296863333333333
$ws.Range("H4").Value = 0.389059
$ws.Range("I4").Value = 0.2137501249891493
$ws.Range("J4").Value = 0.2137501249891493
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.288105
$ws.Range("N4").Value = 0.8643149999999999
$ws.Range("O4").Value = 0.7917647996438357
$ws.Range("P4").Value = 0.7917647996438356
$ws.Range("Q4").Value = 0.037363281065
$ws.Range("R4").Value = 0.3362695295849999
$ws.Range("S4").Value = 0.1692398248858786
$ws.Range("T4").Value = 0.1692398248858786

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tnfsf8"
$ws.Range("C5").Value = "Tnfrsf8"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1296863333333333
$ws.Range("H5").Value = 0.389059
$ws.Range("I5").Value = 0.2137501249891493
$ws.Range("J5").Value = 0.2137501249891493
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.07577199999999999
$ws.Range("N5").Value = 0.227316
$ws.Range("O5").Value = 0.2082352003561643
$ws.Range("P5").Value = 0.2082352003561643
$ws.Range("Q5").Value = 0.009826592849333331
$ws.Range("R5").Value = 0.088439335644
$ws.Range("S5").Value = 0.04451030010327066
$ws.Range("T5").Value = 0.04451030010327067

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Tnfsf8"
$ws.Range("C6").Value = "Tnfrsf8"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.017113
$ws.Range("H6").Value = 0.051339
$ws.Range("I6").Value = 0.02820579312345412
$ws.Range("J6").Value = 0.02820579312345413
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.288105
$ws.Range("N6").Value = 0.8643149999999999
$ws.Range("O6").Value = 0.7917647996438357
$ws.Range("P6").Value = 0.7917647996438356
$ws.Range("Q6").Value = 0.004930340865
$ws.Range("R6").Value = 0.044373067785
$ws.Range("S6").Value = 0.02233235414118713
$ws.Range("T6").Value = 0.02233235414118713

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Tnfsf8"
$ws.Range("C7").Value = "Tnfrsf8"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.017113
$ws.Range("H7").Value = 0.051339
$ws.Range("I7").Value = 0.02820579312345412
$ws.Range("J7").Value = 0.02820579312345413
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.07577199999999999
$ws.Range("N7").Value = 0.227316
$ws.Range("O7").Value = 0.2082352003561643
$ws.Range("P7").Value = 0.2082352003561643
$ws.Range("Q7").Value = 0.001296686236
$ws.Range("R7").Value = 0.011670176124
$ws.Range("S7").Value = 0.005873438982266989
$ws.Range("T7").Value = 0.005873438982266991

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Tnfsf8"
$ws.Range("C8").Value = "Tnfrsf8"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.3848009999999999
$ws.Range("H8").Value = 1.154403
$ws.Range("I8").Value = 0.6342323029099671
$ws.Range("J8").Value = 0.6342323029099671
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.288105
$ws.Range("N8").Value = 0.8643149999999999
$ws.Range("O8").Value = 0.7917647996438357
$ws.Range("P8").Value = 0.7917647996438356
$ws.Range("Q8").Value = 0.110863092105
$ws.Range("R8").Value = 0.9977678289449998
$ws.Range("S8").Value = 0.5021628122411587
$ws.Range("T8").Value = 0.5021628122411586

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Tnfsf8"
$ws.Range("C9").Value = "Tnfrsf8"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.3848009999999999
$ws.Range("H9").Value = 1.154403
$ws.Range("I9").Value = 0.6342323029099671
$ws.Range("J9").Value = 0.6342323029099671
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.07577199999999999
$ws.Range("N9").Value = 0.227316
$ws.Range("O9").Value = 0.2082352003561643
$ws.Range("P9").Value = 0.2082352003561643
$ws.Range("Q9").Value = 0.02915714137199999
$ws.Range("R9").Value = 0.262414272348
$ws.Range("S9").Value = 0.1320694906688085
$ws.Range("T9").Value = 0.1320694906688085
